$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cached value updates from "incentive based DR" re-run
$updates = @(
    @{ Row = 25; Col = 5; Value = 15.6 }
    @{ Row = 25; Col = 6; Value = 16 }
    @{ Row = 25; Col = 11; Value = 64.73416666666667 }
    @{ Row = 26; Col = 6; Value = 2.051282051282053 }
    @{ Row = 26; Col = 11; Value = 49.61628205128206 }
    @{ Row = 95; Col = 8; Value = 48.57500000000001 }
    @{ Row = 95; Col = 10; Value = 8.774999999999999 }
    @{ Row = 95; Col = 11; Value = 47.92041666666667 }
    @{ Row = 96; Col = 8; Value = 48.57500000000001 }
    @{ Row = 97; Col = 8; Value = 48.57500000000001 }
    @{ Row = 98; Col = 8; Value = 48.57500000000001 }
    @{ Row = 99; Col = 8; Value = 48.57500000000001 }
    @{ Row = 100; Col = 8; Value = 76.42500000000001 }
    @{ Row = 101; Col = 8; Value = 99.87500000000001 }
    @{ Row = 101; Col = 9; Value = 23.45 }
    @{ Row = 101; Col = 11; Value = 71.67416666666668 }
    @{ Row = 102; Col = 8; Value = 99.87500000000001 }
    @{ Row = 103; Col = 8; Value = 99.87500000000001 }
    @{ Row = 104; Col = 8; Value = 99.62500000000001 }
    @{ Row = 105; Col = 8; Value = 99.30000000000001 }
    @{ Row = 106; Col = 8; Value = 98.05000000000001 }
    @{ Row = 107; Col = 8; Value = 96.37500000000001 }
    @{ Row = 108; Col = 8; Value = 94.35000000000001 }
    @{ Row = 109; Col = 5; Value = 30.76923076923077 }
    @{ Row = 109; Col = 7; Value = 16 }
    @{ Row = 109; Col = 8; Value = 91.97500000000001 }
    @{ Row = 109; Col = 11; Value = 14.84583333333333 }
    @{ Row = 110; Col = 5; Value = 30.76923076923077 }
    @{ Row = 110; Col = 8; Value = 89.02500000000001 }
    @{ Row = 111; Col = 5; Value = 33.63076923076923 }
    @{ Row = 111; Col = 6; Value = 2.934911242603551 }
    @{ Row = 111; Col = 8; Value = 86 }
    @{ Row = 111; Col = 11; Value = 30.69199457593688 }
    @{ Row = 112; Col = 5; Value = 33.63076923076923 }
    @{ Row = 112; Col = 6; Value = 0 }
    @{ Row = 112; Col = 8; Value = 82.075 }
    @{ Row = 112; Col = 11; Value = 35.6675 }
    @{ Row = 113; Col = 5; Value = 17.22051282051282 }
    @{ Row = 113; Col = 8; Value = 76.90000000000001 }
    @{ Row = 114; Col = 5; Value = 0.8102564102564127 }
    @{ Row = 114; Col = 8; Value = 68.125 }
    @{ Row = 115; Col = 6; Value = 16 }
    @{ Row = 115; Col = 7; Value = 0 }
    @{ Row = 115; Col = 8; Value = 86.375 }
    @{ Row = 115; Col = 9; Value = 18.25 }
    @{ Row = 115; Col = 10; Value = 0 }
    @{ Row = 115; Col = 11; Value = 0 }
    @{ Row = 116; Col = 8; Value = 78.47499999999999 }
    @{ Row = 117; Col = 8; Value = 70.45 }
    @{ Row = 118; Col = 8; Value = 61.275 }
    @{ Row = 119; Col = 8; Value = 52.5 }
    @{ Row = 120; Col = 8; Value = 52.5 }
    @{ Row = 121; Col = 8; Value = 52.5 }
    @{ Row = 122; Col = 8; Value = 52.5 }
    @{ Row = 123; Col = 8; Value = 52.5 }
    @{ Row = 124; Col = 8; Value = 52.5 }
    @{ Row = 124; Col = 9; Value = 0 }
    @{ Row = 124; Col = 11; Value = 45.46666666666667 }
    @{ Row = 125; Col = 8; Value = 52.5 }
    @{ Row = 126; Col = 8; Value = 52.5 }
    @{ Row = 127; Col = 8; Value = 52.5 }
    @{ Row = 128; Col = 8; Value = 52.5 }
    @{ Row = 129; Col = 8; Value = 52.5 }
    @{ Row = 130; Col = 8; Value = 52.5 }
    @{ Row = 131; Col = 8; Value = 52.5 }
    @{ Row = 132; Col = 8; Value = 51.15 }
    @{ Row = 133; Col = 8; Value = 49.125 }
    @{ Row = 134; Col = 8; Value = 46.375 }
    @{ Row = 135; Col = 8; Value = 42.4 }
    @{ Row = 136; Col = 10; Value = 4.55 }
    @{ Row = 136; Col = 11; Value = 37.06666666666668 }
    @{ Row = 366; Col = 8; Value = 65.55 }
    @{ Row = 366; Col = 9; Value = 7.974999999999994 }
    @{ Row = 366; Col = 11; Value = 64.58291666666666 }
    @{ Row = 367; Col = 8; Value = 65.55 }
    @{ Row = 368; Col = 10; Value = 0.25 }
    @{ Row = 368; Col = 11; Value = 46.4125 }
    @{ Row = 435; Col = 8; Value = 2.449999999999985 }
    @{ Row = 435; Col = 9; Value = 2.449999999999985 }
    @{ Row = 435; Col = 11; Value = 46.95583333333331 }
    @{ Row = 436; Col = 8; Value = 30.29999999999999 }
    @{ Row = 437; Col = 8; Value = 58.39999999999998 }
    @{ Row = 438; Col = 8; Value = 86.42499999999998 }
    @{ Row = 439; Col = 8; Value = 86.42499999999998 }
    @{ Row = 440; Col = 8; Value = 86.42499999999998 }
    @{ Row = 441; Col = 8; Value = 86.42499999999998 }
    @{ Row = 442; Col = 10; Value = 1.25 }
    @{ Row = 442; Col = 11; Value = 44.56541666666667 }
    @{ Row = 651; Col = 8; Value = 26.64750000000246 }
    @{ Row = 651; Col = 9; Value = 2.447500000002464 }
    @{ Row = 651; Col = 11; Value = 47.82833333333579 }
    @{ Row = 652; Col = 8; Value = 26.64750000000246 }
    @{ Row = 653; Col = 8; Value = 26.64750000000246 }
    @{ Row = 654; Col = 8; Value = 26.64750000000246 }
    @{ Row = 655; Col = 8; Value = 26.64750000000246 }
    @{ Row = 656; Col = 8; Value = 25.89750000000246 }
    @{ Row = 657; Col = 8; Value = 25.39750000000246 }
    @{ Row = 658; Col = 8; Value = 24.57250000000246 }
    @{ Row = 659; Col = 8; Value = 23.32250000000246 }
    @{ Row = 660; Col = 8; Value = 21.97250000000246 }
    @{ Row = 661; Col = 8; Value = 19.94750000000246 }
    @{ Row = 662; Col = 8; Value = 17.19750000000246 }
    @{ Row = 663; Col = 8; Value = 17.19750000000246 }
    @{ Row = 664; Col = 8; Value = 17.19750000000246 }
    @{ Row = 665; Col = 8; Value = 17.19750000000246 }
    @{ Row = 666; Col = 8; Value = 10.74750000000246 }
    @{ Row = 667; Col = 8; Value = 3.322500000002464 }
    @{ Row = 668; Col = 10; Value = 3.322500000002464 }
    @{ Row = 668; Col = 11; Value = 52.96749999999754 }
    @{ Row = 674; Col = 8; Value = 49.575 }
    @{ Row = 674; Col = 9; Value = 23.3 }
    @{ Row = 674; Col = 11; Value = 66.04708333333333 }
    @{ Row = 675; Col = 9; Value = 27.675 }
    @{ Row = 675; Col = 11; Value = 72.58083333333333 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

Write-Host "Applied" $updates.Count "cell updates"